# Updated symbol list on Tue Dec 20 22:41:09 UTC 2022 with GitHub Actions
# Applies the refreshed price/volume-label values to the "cryptos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates -------------------------------------------
# Values are stored as text in the workbook, so each one is written with a
# leading apostrophe to keep Excel from re-interpreting it as a number
# (which would drop meaningful trailing zeros, e.g. "250.80").
$ws.Range("D2").Value  = "'250.80"
$ws.Range("D3").Value  = "'22.77"
$ws.Range("D5").Value  = "'0.05671"
$ws.Range("D6").Value  = "'3.413"
$ws.Range("D7").Value  = "'6.380"
$ws.Range("D9").Value  = "'0.9280"
$ws.Range("D10").Value = "'0.1440"
$ws.Range("D11").Value = "'0.07433"
$ws.Range("D12").Value = "'0.03186"
$ws.Range("D13").Value = "'0.03073"
$ws.Range("D15").Value = "'3.729"
$ws.Range("D16").Value = "'0.001597"
$ws.Range("D17").Value = "'0.04762"
$ws.Range("D18").Value = "'0.0005789"
$ws.Range("D23").Value = "'3.708"
$ws.Range("D24").Value = "'2.180"
$ws.Range("D25").Value = "'0.3304"
$ws.Range("D26").Value = "'0.1308"
$ws.Range("D28").Value = "'0.0003000"
$ws.Range("D40").Value = "'0.04020"
$ws.Range("D44").Value = "'0.007558"
$ws.Range("D45").Value = "'0.00005802"
$ws.Range("D47").Value = "'0.5000"

# --- Volume(1h) label column (E) updates ---------------------------------
# "Worstin24h" / "Bestin24h" suffix moved from the One (row 18) coin to the
# CoinbaseStockToken (row 47) coin.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
